# "points table index removed"
# The PLAYER NUMBER column (the points-table row index column, column G)
# is removed entirely, and the row that used to read "ranjan ramanayake"
# is corrected to the player's actual name, "Rassie Van der Dussen".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire "PLAYER NUMBER" index column (column G).
$ws.Columns.Item(7).Delete()

# Fix the mis-entered player name in row 4.
$ws.Range("A4").Value = "Rassie Van der Dussen"
